$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.302.02"
$ws.Range("E2").Value = "  +0.31%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.869.25"
$ws.Range("E3").Value = "  +0.41%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  +0.10%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "234.95"
$ws.Range("E5").Value = "  -0.54%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  +0.13%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4698"
$ws.Range("E7").Value = "  +0.44%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2864"
$ws.Range("E8").Value = "  +0.55%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06569"

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.58"
$ws.Range("E10").Value = "  -0.95%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08022"
$ws.Range("E11").Value = "  +1.57%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "96.91"
$ws.Range("E12").Value = "  -0.44%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.869.68"
$ws.Range("E13").Value = "  +0.44%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.112"
$ws.Range("E14").Value = "  -0.96%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6824"
$ws.Range("E15").Value = "  +0.40%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "268.55"
$ws.Range("E16").Value = "  -3.76%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "30.317.49"
$ws.Range("E17").Value = "  +0.38%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.98"
$ws.Range("E18").Value = "  +3.57%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007620"
$ws.Range("E19").Value = "  +4.23%  "

$ws.Range("E20").Value = "  +0.12%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.116.18"
$ws.Range("E21").Value = "  +0.35%  "

$ws.Range("E22").Value = "  +0.20%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.268"
$ws.Range("E23").Value = "  -1.68%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.200"
$ws.Range("E24").Value = "  +0.64%  "

$ws.Range("B25").Value = "Monero"
$ws.Range("C25").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "168.72"
$ws.Range("E25").Value = "  +0.31%  "

$ws.Range("B26").Value = "Cosmos"
$ws.Range("C26").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.378"
$ws.Range("E26").Value = "  +1.56%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.85"
$ws.Range("E27").Value = "  -1.11%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.944"
$ws.Range("E28").Value = "  +0.62%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.369"
$ws.Range("E29").Value = "  -0.74%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.09894"
$ws.Range("E30").Value = "  +1.60%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.368"
$ws.Range("E31").Value = "  +0.04%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.459"
$ws.Range("E32").Value = "  -1.32%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.056"
$ws.Range("E33").Value = "  +0.18%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04682"
$ws.Range("E34").Value = "  -0.91%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.134"
$ws.Range("E35").Value = "  -0.01%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6994"

$ws.Range("E37").Value = "  -0.12%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01867"
$ws.Range("E38").Value = "  +0.11%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.626"
$ws.Range("E39").Value = "  +0.13%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.270"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "72.03"
$ws.Range("E41").Value = "  -3.44%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.952"
$ws.Range("E42").Value = "  -0.02%  "

$ws.Range("B43").Value = "PaxDollar"
$ws.Range("C43").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.001"
$ws.Range("E43").Value = "  +0.13%  "

$ws.Range("B44").Value = "TheSandbox"
$ws.Range("C44").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.4160"
$ws.Range("E44").Value = "  -0.37%  "

$ws.Range("B45").Value = "TrustWalletToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.8397"
$ws.Range("E45").Value = "  -1.02%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "102.72"

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.171"

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.046"
$ws.Range("E48").Value = "  -1.99%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "910.13"
$ws.Range("E49").Value = "  -6.07%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "34.40"
$ws.Range("E50").Value = "  +0.86%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05677"
$ws.Range("E51").Value = "  +0.67%  "
